$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.07"
$ws.Range("E2").Value = "'8.14%"
$ws.Range("D3").Value = "'48.37"
$ws.Range("E3").Value = "'14.95%"
$ws.Range("D4").Value = "'5.256"
$ws.Range("E4").Value = "'4.95%"
$ws.Range("D5").Value = "'0.08103"
$ws.Range("E5").Value = "'7.71%"
$ws.Range("D6").Value = "'4.589"
$ws.Range("E6").Value = "'5.07%"
$ws.Range("D7").Value = "'1.642"
$ws.Range("E7").Value = "'2.81%"
$ws.Range("D8").Value = "'1.200"
$ws.Range("E8").Value = "'30.30%"
$ws.Range("D9").Value = "'0.1300"
$ws.Range("E9").Value = "'9.89%"
$ws.Range("D10").Value = "'0.1947"
$ws.Range("E10").Value = "'6.39%"
$ws.Range("D11").Value = "'0.09517"
$ws.Range("E11").Value = "'6.46%"
$ws.Range("D12").Value = "'0.04617"
$ws.Range("E12").Value = "'11.89%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'-0.16%"
$ws.Range("D14").Value = "'0.001326"
$ws.Range("E14").Value = "'3.41%"
$ws.Range("D15").Value = "'0.005927"
$ws.Range("E15").Value = "'0.08%"
$ws.Range("D16").Value = "'3.343"
$ws.Range("E16").Value = "'0.02%"
$ws.Range("E17").Value = "'1.15%"
$ws.Range("D18").Value = "'0.3403"
$ws.Range("E18").Value = "'2.17%"
$ws.Range("D19").Value = "'8.117"
$ws.Range("E19").Value = "'-2.31%"
$ws.Range("E20").Value = "'4.30%"
$ws.Range("D21").Value = "'0.3126"
$ws.Range("E21").Value = "'0.72%"
$ws.Range("D22").Value = "'0.04289"
$ws.Range("E22").Value = "'4.64%"
$ws.Range("D23").Value = "'0.001310"
$ws.Range("E23").Value = "'3.36%"
$ws.Range("D24").Value = "'0.004250"
$ws.Range("E24").Value = "'9.33%"
$ws.Range("E25").Value = "'3.74%"
$ws.Range("D26").Value = "'0.0003541"
$ws.Range("E26").Value = "'-4.92%"
$ws.Range("D38").Value = "'0.02669"
$ws.Range("E38").Value = "'11.47%"
$ws.Range("D39").Value = "'0.05622"
$ws.Range("E39").Value = "'7.64%"
$ws.Range("D40").Value = "'0.006303"
$ws.Range("E40").Value = "'-9.56%"
$ws.Range("D41").Value = "'0.007690"
$ws.Range("E41").Value = "'-1.07%"
$ws.Range("D42").Value = "'0.1440"
$ws.Range("E42").Value = "'8.71%"
$ws.Range("D43").Value = "'0.007697"
$ws.Range("E43").Value = "'3.89%"
$ws.Range("E44").Value = "'13.77%"
$ws.Range("E45").Value = "'7.12%"
$ws.Range("D46").Value = "'0.00006994"
$ws.Range("E46").Value = "'6.22%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'0.05344"
$ws.Range("E48").Value = "'17.75%"
$ws.Range("D49").Value = "'0.004002"
$ws.Range("E49").Value = "'-4.84%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.08%"
